$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.543145160119082
$ws.Range("C2").Value = 0.6045642687506074
$ws.Range("D2").Value = 0.05575800697469191
$ws.Range("E2").Value = 0.09227724228551359
$ws.Range("F2").Value = 2.705411538345544
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 1.262895537521509
$ws.Range("J2").Value = 0.1479057805073616
$ws.Range("M2").Value = 0.5724660458667543
$ws.Range("N2").Value = 1.799457959988786

$ws.Range("B3").Value = 1.440653581443883
$ws.Range("C3").Value = 0.565724447271748
$ws.Range("D3").Value = 0.05560136927945081
$ws.Range("E3").Value = 0.0923741770050075
$ws.Range("F3").Value = 2.681179437664213
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 1.258904852297867
$ws.Range("J3").Value = 0.148147194015344
$ws.Range("M3").Value = 0.5481573466469456
$ws.Range("N3").Value = 1.819625327823573

$ws.Range("B4").Value = 1.378553937137156
$ws.Range("C4").Value = 0.5422353028174598
$ws.Range("D4").Value = 0.05550901246319384
$ws.Range("E4").Value = 0.09246258146399988
$ws.Range("F4").Value = 2.66792642449731
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 1.257205323389144
$ws.Range("J4").Value = 0.1483581731651604
$ws.Range("M4").Value = 0.5335595744811528
$ws.Range("N4").Value = 1.832670823019775

$ws.Range("B5").Value = 1.353455744878147
$ws.Range("C5").Value = 0.5327527696101697
$ws.Range("D5").Value = 0.05547234174530224
$ws.Range("E5").Value = 0.09250588172466578
$ws.Range("F5").Value = 2.662933025250439
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 1.256700839107523
$ws.Range("J5").Value = 0.1484599138987512
$ws.Range("M5").Value = 0.527693066422529
$ws.Range("N5").Value = 1.838153240818649

$ws.Range("B6").Value = 1.34930073989301
$ws.Range("C6").Value = 0.5311835864488614
$ws.Range("D6").Value = 0.05546631106280842
$ws.Range("E6").Value = 0.09251351140551911
$ws.Range("F6").Value = 2.66212843345734
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 1.256628409596388
$ws.Range("J6").Value = 0.1484777596384177
$ws.Range("M6").Value = 0.5267238951603375
$ws.Range("N6").Value = 1.83907362424122

$ws.Range("B7").Value = 1.378214613528428
$ws.Range("C7").Value = 0.5421070566614219
$ws.Range("D7").Value = 0.05550851399221202
$ws.Range("E7").Value = 0.09246313595635058
$ws.Range("F7").Value = 2.66785743452489
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 1.257197759114923
$ws.Range("J7").Value = 0.1483594814622897
$ws.Range("M7").Value = 0.5334801241830291
$ws.Range("N7").Value = 1.832744088081487

$ws.Range("B8").Value = 1.507633062785487
$ws.Range("C8").Value = 0.5910974175969841
$ws.Range("D8").Value = 0.05570320845572496
$ws.Range("E8").Value = 0.09230467650261609
$ws.Range("F8").Value = 2.696718020902367
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 1.261363300069355
$ws.Range("J8").Value = 0.1479759872091222
$ws.Range("M8").Value = 0.5640162006656908
$ws.Range("N8").Value = 1.806273922297621

$ws.Range("B9").Value = 1.768070890135675
$ws.Range("C9").Value = 0.6900508364707889
$ws.Range("D9").Value = 0.05611512934374829
$ws.Range("E9").Value = 0.09222271856961406
$ws.Range("F9").Value = 2.766284719212024
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 1.275522736755136
$ws.Range("J9").Value = 0.1477226294619527
$ws.Range("M9").Value = 0.6265138492158542
$ws.Range("N9").Value = 1.759639107451903

$ws.Range("B10").Value = 1.963570924416956
$ws.Range("C10").Value = 0.7645703173393485
$ws.Range("D10").Value = 0.05643595357482489
$ws.Range("E10").Value = 0.09230151436973522
$ws.Range("F10").Value = 2.825413979438821
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 1.289627965807
$ws.Range("J10").Value = 0.1478417792982327
$ws.Range("M10").Value = 0.6740526506733602
$ws.Range("N10").Value = 1.728611143821013

$ws.Range("B11").Value = 2.053435510818019
$ws.Range("C11").Value = 0.7988797615585668
$ws.Range("D11").Value = 0.05658582421271241
$ws.Range("E11").Value = 0.09236745904882504
$ws.Range("F11").Value = 2.854079502260191
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 1.296859878005179
$ws.Range("J11").Value = 0.1479625768201487
$ws.Range("M11").Value = 0.6960376782652986
$ws.Range("N11").Value = 1.715202477115582

$ws.Range("B12").Value = 2.087600248394892
$ws.Range("C12").Value = 0.8119317994283506
$ws.Range("D12").Value = 0.05664313792930287
$ws.Range("E12").Value = 0.09239675067127529
$ws.Range("F12").Value = 2.865190407876611
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 1.299716532482719
$ws.Range("J12").Value = 0.1480179187246904
$ws.Range("M12").Value = 0.704414916922417
$ws.Range("N12").Value = 1.710226916638547

$ws.Range("B13").Value = 2.080236239506689
$ws.Range("C13").Value = 0.8091181369487686
$ws.Range("D13").Value = 0.05663076950060209
$ws.Range("E13").Value = 0.0923902502052254
$ws.Range("F13").Value = 2.86278606633897
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 1.299096035459627
$ws.Range("J13").Value = 0.1480055725994731
$ws.Range("M13").Value = 0.702608411625917
$ws.Range("N13").Value = 1.711293948037358

$ws.Range("B14").Value = 2.056243552686794
$ws.Range("C14").Value = 0.7999523574416685
$ws.Range("D14").Value = 0.056590528223742
$ws.Range("E14").Value = 0.09236978236098814
$ws.Range("F14").Value = 2.854988464500764
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 1.297092525230767
$ws.Range("J14").Value = 0.1479669373343313
$ws.Range("M14").Value = 0.6967258350692163
$ws.Range("N14").Value = 1.714791087578973

$ws.Range("B15").Value = 2.041564950718509
$ws.Range("C15").Value = 0.7943458645093529
$ws.Range("D15").Value = 0.05656595223157623
$ws.Range("E15").Value = 0.09235780753241407
$ws.Range("F15").Value = 2.85024558576265
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 1.295880720710855
$ws.Range("J15").Value = 0.1479445227965286
$ws.Range("M15").Value = 0.6931293664170823
$ws.Range("N15").Value = 1.716946485492137

$ws.Range("B16").Value = 1.957716886245862
$ws.Range("C16").Value = 0.7623364372832953
$ws.Range("D16").Value = 0.05642623786400236
$ws.Range("E16").Value = 0.09229780968316348
$ws.Range("F16").Value = 2.823576332608724
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 1.289171821980219
$ws.Range("J16").Value = 0.1478352268465031
$ws.Range("M16").Value = 0.6726231351886724
$ws.Range("N16").Value = 1.729501680687221

$ws.Range("B17").Value = 1.906517977226088
$ws.Range("C17").Value = 0.7428053174047022
$ws.Range("D17").Value = 0.05634153052908175
$ws.Range("E17").Value = 0.09226870612542548
$ws.Range("F17").Value = 2.807669504100261
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 1.285265558106985
$ws.Range("J17").Value = 0.1477852495935608
$ws.Range("M17").Value = 0.6601354893806928
$ws.Range("N17").Value = 1.73738507168401

$ws.Range("B18").Value = 1.877157244714056
$ws.Range("C18").Value = 0.7316100901893492
$ws.Range("D18").Value = 0.05629317884962504
$ws.Range("E18").Value = 0.09225480018508314
$ws.Range("F18").Value = 2.798686559564374
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 1.283095474487169
$ws.Range("J18").Value = 0.1477627709346407
$ws.Range("M18").Value = 0.6529867342390645
$ws.Range("N18").Value = 1.741985823285212

$ws.Range("B19").Value = 1.867231203497795
$ws.Range("C19").Value = 0.7278261761078966
$ws.Range("D19").Value = 0.05627687141957693
$ws.Range("E19").Value = 0.09225057882082588
$ws.Range("F19").Value = 2.795673588156546
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 1.282373868761439
$ws.Range("J19").Value = 0.1477562357357485
$ws.Range("M19").Value = 0.6505720860365045
$ws.Range("N19").Value = 1.743554954474568

$ws.Range("B20").Value = 1.911959121727477
$ws.Range("C20").Value = 0.7448804421615591
$ws.Range("D20").Value = 0.05635050952986909
$ws.Range("E20").Value = 0.0922715110392236
$ws.Range("F20").Value = 2.809345595062553
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 1.285673442793879
$ws.Range("J20").Value = 0.1477899209954643
$ws.Range("M20").Value = 0.6614613193624592
$ws.Range("N20").Value = 1.736538992232219

$ws.Range("B21").Value = 2.063287114089405
$ws.Range("C21").Value = 0.8026429396424533
$ws.Range("D21").Value = 0.05660233286550209
$ws.Range("E21").Value = 0.09237567709154781
$ws.Range("F21").Value = 2.857271850155485
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 1.297677793196911
$ws.Range("J21").Value = 0.1479780247680793
$ws.Range("M21").Value = 0.6984522770304125
$ws.Range("N21").Value = 1.713761120573537

$ws.Range("B22").Value = 2.162975963051792
$ws.Range("C22").Value = 0.8407429971734928
$ws.Range("D22").Value = 0.05677018220914931
$ws.Range("E22").Value = 0.09246893097948039
$ws.Range("F22").Value = 2.890086810185494
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 1.306212023884939
$ws.Range("J22").Value = 0.1481569196572039
$ws.Range("M22").Value = 0.7229311685456423
$ws.Range("N22").Value = 1.699469285952233

$ws.Range("B23").Value = 2.109697727189882
$ws.Range("C23").Value = 0.8203760790838146
$ws.Range("D23").Value = 0.05668029992596857
$ws.Range("E23").Value = 0.0924168588307559
$ws.Range("F23").Value = 2.872435708274963
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 1.30159384935402
$ws.Range("J23").Value = 0.1480563124134449
$ws.Range("M23").Value = 0.7098384803820181
$ws.Range("N23").Value = 1.707042522424757

$ws.Range("B24").Value = 1.909498949260808
$ws.Range("C24").Value = 0.7439421740572243
$ws.Range("D24").Value = 0.05634644904063713
$ws.Range("E24").Value = 0.092270234134471
$ws.Range("F24").Value = 2.808587329555877
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 1.285488802465267
$ws.Range("J24").Value = 0.1477877895769453
$ws.Range("M24").Value = 0.6608618164668059
$ws.Range("N24").Value = 1.736921291472257

$ws.Range("B25").Value = 1.696893226698364
$ws.Range("C25").Value = 0.6629664119413974
$ws.Range("D25").Value = 0.05600048829814774
$ws.Range("E25").Value = 0.09222044641369465
$ws.Range("F25").Value = 2.746065047833255
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 1.271045733094752
$ws.Range("J25").Value = 0.1477376524002736
$ws.Range("M25").Value = 0.6093237394563076
$ws.Range("N25").Value = 1.771688500391377

